# Add a new worksheet "ODI Batting Extra" to the workbook, as the last
# sheet (right after "ODI Bowling"), and populate it with a header row plus
# one data row, matching the PlayerPerformance_4746 export schema.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 data. BATTING_POSITION is numeric; the rest are text (including the
# "0.99%" percent-looking string, which must stay literal text, not get
# auto-converted to a numeric percentage by Excel's smart entry).
$newSheet.Cells.Item(2, 1).Value = "4238"
$newSheet.Cells.Item(2, 2).Value = 7
$newSheet.Cells.Item(2, 3).Value = "0"
$newSheet.Cells.Item(2, 4).Value = "0"

$pctCell = $newSheet.Cells.Item(2, 5)
$pctCell.NumberFormat = "@"
$pctCell.Value = "0.99%"
$pctCell.Style = "Normal"

$newSheet.Cells.Item(2, 6).Value = "NO"
